# Actualización 10 de Mayo
# Update statistics across the "Estadisticos 1P", "Estadisticos 2P",
# "Estadisticos Final" sheets, and add a new "rescatable" student row
# to the "Rescatables" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Estadisticos 1P
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

# Row 3 (4ARHV)
$ws1.Range("D3").Value = 12
$ws1.Range("F3").Value = 16
$ws1.Range("G3").Value = 55.17

# Row 7 (4ASV)
$ws1.Range("D7").Value = 5
$ws1.Range("F7").Value = 28
$ws1.Range("G7").Value = 82.34999999999999

# ---------------------------------------------------------------
# Sheet: Estadisticos 2P
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

# Row 2 (2ARHV)
$ws2.Range("D2").Value = 11
$ws2.Range("E2").Value = 9
$ws2.Range("F2").Value = 29
$ws2.Range("G2").Value = 72.5
$ws2.Range("H2").Value = 7.4

# Row 3 (4AEV)
$ws2.Range("D3").Value = 17
$ws2.Range("E3").Value = 6
$ws2.Range("F3").Value = 11
$ws2.Range("G3").Value = 37.93
$ws2.Range("H3").Value = 6.3

# Row 4 (4ALCV)
$ws2.Range("D4").Value = 8
$ws2.Range("E4").Value = 7
$ws2.Range("F4").Value = 25
$ws2.Range("G4").Value = 75.76000000000001

# Row 5 (4APV)
$ws2.Range("D5").Value = 19
$ws2.Range("E5").Value = 13
$ws2.Range("F5").Value = 19
$ws2.Range("G5").Value = 50
$ws2.Range("H5").Value = 7.6

# Row 6 (4ARHV)
$ws2.Range("D6").Value = 14
$ws2.Range("E6").Value = 7
$ws2.Range("F6").Value = 25
$ws2.Range("G6").Value = 64.09999999999999

# Row 7 (4ASV)
$ws2.Range("D7").Value = 10
$ws2.Range("E7").Value = 5
$ws2.Range("F7").Value = 24
$ws2.Range("G7").Value = 70.59
$ws2.Range("H7").Value = 7

# ---------------------------------------------------------------
# Sheet: Estadisticos Final
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

# Row 3 (4ARHV)
$ws3.Range("D3").Value = 12
$ws3.Range("F3").Value = 16
$ws3.Range("G3").Value = 55.17
$ws3.Range("H3").Value = 6.6

# Row 7 (4ASV)
$ws3.Range("D7").Value = 5
$ws3.Range("E7").Value = 0
$ws3.Range("F7").Value = 29
$ws3.Range("G7").Value = 85.29000000000001

# ---------------------------------------------------------------
# Sheet: Rescatables - add a new student row
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$ws4.Cells.Item(2, 1).Value = 19330051920441
$ws4.Cells.Item(2, 2).Value = "GARCIA"
$ws4.Cells.Item(2, 3).Value = "ANTONIO"
$ws4.Cells.Item(2, 4).Value = "ABRAHAM"
$ws4.Cells.Item(2, 5).Value = "INGLÉS IV"
$ws4.Cells.Item(2, 6).Value = "4APV"
$ws4.Cells.Item(2, 7).Value = 2
